$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParaIndex($pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like $pattern) {
            return $i
        }
    }
    return -1
}

function Replace-ParaXml($idx, $xml) {
    $para = $d.Paragraphs($idx)
    $r = $para.Range
    $r2 = $d.Range($r.Start, $r.End - 1)
    $r2.Collapse(0)
    $r2.InsertXML($xml)
    $orig = $d.Paragraphs($idx)
    $orig.Range.Delete()
}

# ---------------------------------------------------------------------------
# Change 4: move <w:lastRenderedPageBreak/> from the "More damage, range, "
# paragraph (Turret upgrades / Machine gun, numId 11, ilvl 3) to the
# "Machine gun" paragraph right above it (ilvl 2).
# ---------------------------------------------------------------------------
$idxDamageRange = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*More damage, range*") {
        $prevT = $d.Paragraphs($i - 1).Range.Text
        if ($prevT -like "*Machine gun*") {
            $idxDamageRange = $i
            break
        }
    }
}
$idxMachineGun = $idxDamageRange - 1

$machineGunXml = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"2`"/><w:numId w:val=`"11`"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Machine gun</w:t></w:r></w:p>"
Replace-ParaXml $idxMachineGun $machineGunXml

# the "More damage, range" paragraph index is unchanged (we replaced the
# paragraph above it 1-for-1), so it is still $idxDamageRange
$moreDamageXml = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"3`"/><w:numId w:val=`"11`"/></w:numPr></w:pPr><w:r><w:t xml:space=`"preserve`">More damage, range, </w:t></w:r><w:r><w:t>lower costs</w:t></w:r></w:p>"
Replace-ParaXml $idxDamageRange $moreDamageXml

# ---------------------------------------------------------------------------
# Change 3: append to "Gives more money faster" and add the new
# "Player Base" + "More hp, shoots faster, more damage" paragraphs.
# ---------------------------------------------------------------------------
$idxGivesMoney = Find-ParaIndex "*Gives more money faster*"

$givesMoneyXml = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"9`"/></w:numPr></w:pPr><w:r><w:t>Gives more money faster</w:t></w:r><w:r><w:t>, more hp,  gives more money each time</w:t></w:r></w:p>"
Replace-ParaXml $idxGivesMoney $givesMoneyXml

$idxGivesMoney = Find-ParaIndex "*Gives more money faster*"
$anchor = $d.Paragraphs($idxGivesMoney).Range
$anchor2 = $d.Range($anchor.Start, $anchor.End - 1)
$anchor2.Collapse(0)
$newParasXml = "<w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"9`"/></w:numPr></w:pPr><w:r><w:t>Player Base</w:t></w:r></w:p><w:p $wns><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"9`"/></w:numPr></w:pPr><w:r><w:t xml:space=`"preserve`">More hp, shoots faster, </w:t></w:r><w:r><w:t>more damage</w:t></w:r></w:p>"
$anchor2.InsertXML($newParasXml)

# ---------------------------------------------------------------------------
# Change 2: remove the "Gets piercing effect" paragraph entirely.
# ---------------------------------------------------------------------------
$idxPiercing = Find-ParaIndex "*Gets piercing effect*"
$d.Paragraphs($idxPiercing).Range.Delete()

# ---------------------------------------------------------------------------
# Change 1: append ", more hp" as a new run to the "Deals more damage,
# attacks faster" paragraph.
# ---------------------------------------------------------------------------
$idxDealsDamage = Find-ParaIndex "*Deals more damage*"
$para = $d.Paragraphs($idxDealsDamage)
$r = $para.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.InsertAfter(", more hp")
